$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156, shifting rows 156:280 down to 157:281
$ws.Rows(156).Insert()

# Populate the newly inserted row 156 with the latest weekly data point
$ws.Cells.Item(156, 1).Value2 = 8
$ws.Cells.Item(156, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(156, 3).Value2 = "Coquimbo"
$ws.Cells.Item(156, 4).Value2 = 44729
$ws.Cells.Item(156, 5).Value2 = 4
$ws.Cells.Item(156, 6).Value2 = 100112012
$ws.Cells.Item(156, 7).Value2 = "Espinaca"
$ws.Cells.Item(156, 8).Value2 = "Sin especificar"
$ws.Cells.Item(156, 9).Value2 = "Primera"
$ws.Cells.Item(156, 10).Value2 = 3100
$ws.Cells.Item(156, 11).Value2 = 500
$ws.Cells.Item(156, 12).Value2 = 600
$ws.Cells.Item(156, 13).Value2 = 550
$ws.Cells.Item(156, 14).Value2 = "$/atado 300 a 500 gramos"
$ws.Cells.Item(156, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(156, 16).Value2 = 1100
$ws.Cells.Item(156, 17).Value2 = 0.5
$ws.Cells.Item(156, 18).Value2 = "Hortaliza"
